$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '2026-02-01'
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = '18:01:10'
$ws.Range("C8").Value = '18:00'
$ws.Range("D8").Value = 'Bathroom'
$ws.Range("E8").Value = 'No Motion'
$ws.Range("F8").Value = 'Inactive'

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '2026-02-01'
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = '18:01:11'
$ws.Range("C9").Value = '18:00'
$ws.Range("D9").Value = 'Bathroom'
$ws.Range("E9").Value = 'No Motion'
$ws.Range("F9").Value = 'Inactive'

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '2026-02-01'
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = '18:01:15'
$ws.Range("C10").Value = '18:00'
$ws.Range("D10").Value = 'Bathroom'
$ws.Range("E10").Value = 'No Motion'
$ws.Range("F10").Value = 'Inactive'

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '2026-02-01'
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = '18:01:20'
$ws.Range("C11").Value = '18:00'
$ws.Range("D11").Value = 'Bathroom'
$ws.Range("E11").Value = 'No Motion'
$ws.Range("F11").Value = 'Inactive'

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '2026-02-01'
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = '18:01:25'
$ws.Range("C12").Value = '18:00'
$ws.Range("D12").Value = 'Bathroom'
$ws.Range("E12").Value = 'No Motion'
$ws.Range("F12").Value = 'Inactive'

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '2026-02-01'
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = '18:01:30'
$ws.Range("C13").Value = '18:00'
$ws.Range("D13").Value = 'Bathroom'
$ws.Range("E13").Value = 'No Motion'
$ws.Range("F13").Value = 'Inactive'

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '2026-02-01'
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = '18:01:35'
$ws.Range("C14").Value = '18:00'
$ws.Range("D14").Value = 'Bathroom'
$ws.Range("E14").Value = 'No Motion'
$ws.Range("F14").Value = 'Inactive'

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '2026-02-01'
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = '18:01:40'
$ws.Range("C15").Value = '18:00'
$ws.Range("D15").Value = 'Bathroom'
$ws.Range("E15").Value = 'No Motion'
$ws.Range("F15").Value = 'Inactive'

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '2026-02-01'
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = '18:01:45'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Bathroom'
$ws.Range("E16").Value = 'No Motion'
$ws.Range("F16").Value = 'Inactive'

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '2026-02-01'
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = '18:01:50'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Bathroom'
$ws.Range("E17").Value = 'No Motion'
$ws.Range("F17").Value = 'Inactive'

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '2026-02-01'
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = '18:01:55'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Bathroom'
$ws.Range("E18").Value = 'No Motion'
$ws.Range("F18").Value = 'Inactive'

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '2026-02-01'
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = '18:02:00'
$ws.Range("C19").Value = '18:00'
$ws.Range("D19").Value = 'Bathroom'
$ws.Range("E19").Value = 'No Motion'
$ws.Range("F19").Value = 'Inactive'

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '2026-02-01'
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = '18:02:05'
$ws.Range("C20").Value = '18:00'
$ws.Range("D20").Value = 'Bathroom'
$ws.Range("E20").Value = 'No Motion'
$ws.Range("F20").Value = 'Inactive'

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '2026-02-01'
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = '18:01:10'
$ws.Range("C7").Value = '18:00'
$ws.Range("D7").Value = 'Bathroom'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '82.4%'
$ws.Range("E7").ClearFormats()
$ws.Range("F7").Value = 'Active'

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '2026-02-01'
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = '18:01:11'
$ws.Range("C8").Value = '18:00'
$ws.Range("D8").Value = 'Bathroom'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '81.3%'
$ws.Range("E8").ClearFormats()
$ws.Range("F8").Value = 'Active'

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '2026-02-01'
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = '18:01:21'
$ws.Range("C9").Value = '18:00'
$ws.Range("D9").Value = 'Bathroom'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '81.3%'
$ws.Range("E9").ClearFormats()
$ws.Range("F9").Value = 'Active'

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '2026-02-01'
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = '18:01:26'
$ws.Range("C10").Value = '18:00'
$ws.Range("D10").Value = 'Bathroom'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '80.9%'
$ws.Range("E10").ClearFormats()
$ws.Range("F10").Value = 'Active'

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '2026-02-01'
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = '18:01:31'
$ws.Range("C11").Value = '18:00'
$ws.Range("D11").Value = 'Bathroom'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '81.1%'
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = 'Active'

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '2026-02-01'
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = '18:01:36'
$ws.Range("C12").Value = '18:00'
$ws.Range("D12").Value = 'Bathroom'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '82.0%'
$ws.Range("E12").ClearFormats()
$ws.Range("F12").Value = 'Active'

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '2026-02-01'
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = '18:01:41'
$ws.Range("C13").Value = '18:00'
$ws.Range("D13").Value = 'Bathroom'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '81.0%'
$ws.Range("E13").ClearFormats()
$ws.Range("F13").Value = 'Active'

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '2026-02-01'
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = '18:01:46'
$ws.Range("C14").Value = '18:00'
$ws.Range("D14").Value = 'Bathroom'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '81.9%'
$ws.Range("E14").ClearFormats()
$ws.Range("F14").Value = 'Active'

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '2026-02-01'
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = '18:01:51'
$ws.Range("C15").Value = '18:00'
$ws.Range("D15").Value = 'Bathroom'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '81.0%'
$ws.Range("E15").ClearFormats()
$ws.Range("F15").Value = 'Active'

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '2026-02-01'
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = '18:01:56'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Bathroom'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '81.9%'
$ws.Range("E16").ClearFormats()
$ws.Range("F16").Value = 'Active'

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '2026-02-01'
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = '18:02:01'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Bathroom'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '80.9%'
$ws.Range("E17").ClearFormats()
$ws.Range("F17").Value = 'Active'

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '2026-02-01'
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = '18:02:06'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Bathroom'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '81.8%'
$ws.Range("E18").ClearFormats()
$ws.Range("F18").Value = 'Active'

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '2026-02-01'
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = '18:01:11'
$ws.Range("C7").Value = '18:00'
$ws.Range("D7").Value = 'Bathroom'
$ws.Range("E7").Value = '28.6C'
$ws.Range("F7").Value = 'Active'

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '2026-02-01'
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = '18:01:11'
$ws.Range("C8").Value = '18:00'
$ws.Range("D8").Value = 'Bathroom'
$ws.Range("E8").Value = '28.6C'
$ws.Range("F8").Value = 'Active'

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '2026-02-01'
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = '18:01:21'
$ws.Range("C9").Value = '18:00'
$ws.Range("D9").Value = 'Bathroom'
$ws.Range("E9").Value = '28.6C'
$ws.Range("F9").Value = 'Active'

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '2026-02-01'
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = '18:01:26'
$ws.Range("C10").Value = '18:00'
$ws.Range("D10").Value = 'Bathroom'
$ws.Range("E10").Value = '28.6C'
$ws.Range("F10").Value = 'Active'

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '2026-02-01'
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = '18:01:31'
$ws.Range("C11").Value = '18:00'
$ws.Range("D11").Value = 'Bathroom'
$ws.Range("E11").Value = '28.6C'
$ws.Range("F11").Value = 'Active'

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '2026-02-01'
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = '18:01:36'
$ws.Range("C12").Value = '18:00'
$ws.Range("D12").Value = 'Bathroom'
$ws.Range("E12").Value = '28.6C'
$ws.Range("F12").Value = 'Active'

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '2026-02-01'
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = '18:01:41'
$ws.Range("C13").Value = '18:00'
$ws.Range("D13").Value = 'Bathroom'
$ws.Range("E13").Value = '28.6C'
$ws.Range("F13").Value = 'Active'

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '2026-02-01'
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = '18:01:46'
$ws.Range("C14").Value = '18:00'
$ws.Range("D14").Value = 'Bathroom'
$ws.Range("E14").Value = '28.6C'
$ws.Range("F14").Value = 'Active'

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '2026-02-01'
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = '18:01:51'
$ws.Range("C15").Value = '18:00'
$ws.Range("D15").Value = 'Bathroom'
$ws.Range("E15").Value = '28.6C'
$ws.Range("F15").Value = 'Active'

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '2026-02-01'
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = '18:01:56'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Bathroom'
$ws.Range("E16").Value = '28.7C'
$ws.Range("F16").Value = 'Active'

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '2026-02-01'
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = '18:02:01'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Bathroom'
$ws.Range("E17").Value = '28.7C'
$ws.Range("F17").Value = 'Active'

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '2026-02-01'
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = '18:02:06'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Bathroom'
$ws.Range("E18").Value = '28.7C'
$ws.Range("F18").Value = 'Active'

$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '2026-02-01'
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = '18:01:15'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Living Room Main Door'
$ws.Range("E17").Value = 'ENTER'
$ws.Range("F17").Value = 'User ENTERED Living Room Main Door'

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '2026-02-01'
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = '18:01:30'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Living Room Main Door'
$ws.Range("E18").Value = 'EXIT'
$ws.Range("F18").Value = 'User EXITED Living Room Main Door'

$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '2026-02-01'
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = '18:01:16'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Living Room Main Door'
$ws.Range("E16").Value = 'Image Captured'
$ws.Range("F16").Value = 'Active'

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '2026-02-01'
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = '18:01:30'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Living Room Main Door'
$ws.Range("E17").Value = 'Image Captured'
$ws.Range("F17").Value = 'Active'
